# Update "想去人数" (want-to-go count) figures for the refreshed data pull.
# Two sheets carry the same underlying rows and both need the same updates:
#   展览 (Exhibitions)  - rows 2,4,5,6 in column F
#   全部类型 (All types) - rows 2,4,5,7 in column F (row 6 here is the concert,
#                          which is unaffected; the exhibition row shifted to 7)

$wb = $excel.ActiveWorkbook

$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 10418
$wsExpo.Range("F4").Value = 59
$wsExpo.Range("F5").Value = 650
$wsExpo.Range("F6").Value = 486

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 10418
$wsAll.Range("F4").Value = 59
$wsAll.Range("F5").Value = 650
$wsAll.Range("F7").Value = 486
